$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Fbn1"
$ws.Cells.Item(2,3).Value = "Itga5"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 7.848425
$ws.Cells.Item(2,8).Value = 23.545275
$ws.Cells.Item(2,9).Value = 0.02436729568045431
$ws.Cells.Item(2,10).Value = 0.02436729568045431
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 28.72417333333333
$ws.Cells.Item(2,14).Value = 86.17251999999999
$ws.Cells.Item(2,15).Value = 0.4233259107972328
$ws.Cells.Item(2,16).Value = 0.4233259107972328
$ws.Cells.Item(2,17).Value = 225.4395200936666
$ws.Cells.Item(2,18).Value = 2028.955680843
$ws.Cells.Item(2,19).Value = 0.01031530763759379
$ws.Cells.Item(2,20).Value = 0.0103153076375938

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Fbn1"
$ws.Cells.Item(3,3).Value = "Itga5"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 7.848425
$ws.Cells.Item(3,8).Value = 23.545275
$ws.Cells.Item(3,9).Value = 0.02436729568045431
$ws.Cells.Item(3,10).Value = 0.02436729568045431
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 30.56986233333333
$ws.Cells.Item(3,14).Value = 91.709587
$ws.Cells.Item(3,15).Value = 0.4505269713084062
$ws.Cells.Item(3,16).Value = 0.4505269713084062
$ws.Cells.Item(3,17).Value = 239.9252717834917
$ws.Cells.Item(3,18).Value = 2159.327446051425
$ws.Cells.Item(3,19).Value = 0.01097812392189149
$ws.Cells.Item(3,20).Value = 0.01097812392189149

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Fbn1"
$ws.Cells.Item(4,3).Value = "Itga5"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 7.848425
$ws.Cells.Item(4,8).Value = 23.545275
$ws.Cells.Item(4,9).Value = 0.02436729568045431
$ws.Cells.Item(4,10).Value = 0.02436729568045431
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 8.559531999999999
$ws.Cells.Item(4,14).Value = 25.678596
$ws.Cells.Item(4,15).Value = 0.126147117894361
$ws.Cells.Item(4,16).Value = 0.126147117894361
$ws.Cells.Item(4,17).Value = 67.1788449371
$ws.Cells.Item(4,18).Value = 604.6096044339
$ws.Cells.Item(4,19).Value = 0.003073864120969022
$ws.Cells.Item(4,20).Value = 0.003073864120969022

$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Fbn1"
$ws.Cells.Item(5,3).Value = "Itga5"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 291.329961
$ws.Cells.Item(5,8).Value = 873.989883
$ws.Cells.Item(5,9).Value = 0.9045029162236017
$ws.Cells.Item(5,10).Value = 0.9045029162236017
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 28.72417333333333
$ws.Cells.Item(5,14).Value = 86.17251999999999
$ws.Cells.Item(5,15).Value = 0.4233259107972328
$ws.Cells.Item(5,16).Value = 0.4233259107972328
$ws.Cells.Item(5,17).Value = 8368.212296957237
$ws.Cells.Item(5,18).Value = 75313.91067261515
$ws.Cells.Item(5,19).Value = 0.3828995208291093
$ws.Cells.Item(5,20).Value = 0.3828995208291093

$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Fbn1"
$ws.Cells.Item(6,3).Value = "Itga5"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 291.329961
$ws.Cells.Item(6,8).Value = 873.989883
$ws.Cells.Item(6,9).Value = 0.9045029162236017
$ws.Cells.Item(6,10).Value = 0.9045029162236017
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 30.56986233333333
$ws.Cells.Item(6,14).Value = 91.709587
$ws.Cells.Item(6,15).Value = 0.4505269713084062
$ws.Cells.Item(6,16).Value = 0.4505269713084062
$ws.Cells.Item(6,17).Value = 8905.916801345367
$ws.Cells.Item(6,18).Value = 80153.25121210831
$ws.Cells.Item(6,19).Value = 0.4075029593858404
$ws.Cells.Item(6,20).Value = 0.4075029593858404

$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Fbn1"
$ws.Cells.Item(7,3).Value = "Itga5"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 291.329961
$ws.Cells.Item(7,8).Value = 873.989883
$ws.Cells.Item(7,9).Value = 0.9045029162236017
$ws.Cells.Item(7,10).Value = 0.9045029162236017
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 8.559531999999999
$ws.Cells.Item(7,14).Value = 25.678596
$ws.Cells.Item(7,15).Value = 0.126147117894361
$ws.Cells.Item(7,16).Value = 0.126147117894361
$ws.Cells.Item(7,17).Value = 2493.648123738251
$ws.Cells.Item(7,18).Value = 22442.83311364426
$ws.Cells.Item(7,19).Value = 0.114100436008652
$ws.Cells.Item(7,20).Value = 0.114100436008652

$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Fbn1"
$ws.Cells.Item(8,3).Value = "Itga5"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 22.91008466666667
$ws.Cells.Item(8,8).Value = 68.730254
$ws.Cells.Item(8,9).Value = 0.07112978809594397
$ws.Cells.Item(8,10).Value = 0.07112978809594397
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 28.72417333333333
$ws.Cells.Item(8,14).Value = 86.17251999999999
$ws.Cells.Item(8,15).Value = 0.4233259107972328
$ws.Cells.Item(8,16).Value = 0.4233259107972328
$ws.Cells.Item(8,17).Value = 658.0732430466754
$ws.Cells.Item(8,18).Value = 5922.659187420079
$ws.Cells.Item(8,19).Value = 0.03011108233052965
$ws.Cells.Item(8,20).Value = 0.03011108233052965

$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Fbn1"
$ws.Cells.Item(9,3).Value = "Itga5"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 22.91008466666667
$ws.Cells.Item(9,8).Value = 68.730254
$ws.Cells.Item(9,9).Value = 0.07112978809594397
$ws.Cells.Item(9,10).Value = 0.07112978809594397
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 30.56986233333333
$ws.Cells.Item(9,14).Value = 91.709587
$ws.Cells.Item(9,15).Value = 0.4505269713084062
$ws.Cells.Item(9,16).Value = 0.4505269713084062
$ws.Cells.Item(9,17).Value = 700.3581343050108
$ws.Cells.Item(9,18).Value = 6303.223208745098
$ws.Cells.Item(9,19).Value = 0.03204588800067436
$ws.Cells.Item(9,20).Value = 0.03204588800067436

$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Fbn1"
$ws.Cells.Item(10,3).Value = "Itga5"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 22.91008466666667
$ws.Cells.Item(10,8).Value = 68.730254
$ws.Cells.Item(10,9).Value = 0.07112978809594397
$ws.Cells.Item(10,10).Value = 0.07112978809594397
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 8.559531999999999
$ws.Cells.Item(10,14).Value = 25.678596
$ws.Cells.Item(10,15).Value = 0.126147117894361
$ws.Cells.Item(10,16).Value = 0.126147117894361
$ws.Cells.Item(10,17).Value = 196.0996028270426
$ws.Cells.Item(10,18).Value = 1764.896425443384
$ws.Cells.Item(10,19).Value = 0.008972817764739958
$ws.Cells.Item(10,20).Value = 0.008972817764739958
